$d = $word.ActiveDocument

# ---------------------------------------------------------------------
# Both headings in this document ("Scope" and "Outside Scope") are
# styled with Heading3 + bold 14pt "Arial Nova Light" runs. The target
# edit renames them to "Included in the Database Scope" and
# "Excluded from the Database Scope" respectively, while keeping the
# existing "Scope" run's formatting/identity intact and splitting the
# text across multiple runs (as shown by the source diff).
#
# Range.InsertXML (with a full <pkg:package>/<w:document>/<w:body>/<w:p>
# wrapper) only splices new run(s) into an *existing* paragraph cleanly
# when the insertion point is the paragraph's very start; anywhere else
# it inserts a sibling paragraph instead. So for each heading we clear
# the paragraph's text first and then InsertXML the full, correctly
# split set of runs at the (now empty) paragraph's start.
# ---------------------------------------------------------------------

$rPr = '<w:rPr><w:rFonts w:ascii="Arial Nova Light" w:hAnsi="Arial Nova Light"/><w:b/><w:bCs/><w:color w:val="auto"/><w:sz w:val="28"/><w:szCs w:val="28"/></w:rPr>'

function Insert-RunsAtParagraphStart($para, [string]$runsXml) {
    $startPos = $para.Range.Start
    $insPoint = $d.Range($startPos, $startPos)
    $xml = '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p>' + $runsXml + '</w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
    $null = $insPoint.InsertXML($xml)
}

# -----------------------------------------------------------------
# Edit 1: "Scope" heading -> "Included in the Database " + "Scope"
#   The original run (w:rsidRPr="00951ACA") keeps its rPr and becomes
#   the "Included in the Database " run; a brand-new plain run
#   carrying the original "Scope" text is appended after it.
# -----------------------------------------------------------------
$find1 = $d.Content
$find1.Find.ClearFormatting()
$found1 = $find1.Find.Execute("Scope", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if ($found1) {
    $p1 = $find1.Paragraphs.Item(1)

    $p1Full = $d.Range($p1.Range.Start, $p1.Range.End - 1)
    $p1Full.Text = ""

    $runs1 = '<w:r w:rsidRPr="00951ACA">' + $rPr + '<w:t xml:space="preserve">Included in the Database </w:t></w:r>' + `
             '<w:r>' + $rPr + '<w:t>Scope</w:t></w:r>'
    Insert-RunsAtParagraphStart $p1 $runs1
}

# -----------------------------------------------------------------
# Edit 2: "Outside Scope" heading -> "Excluded from the Database" +
#   " " + "Scope"
#   The original first run (plain <w:r>, text "Outside ") becomes
#   "Excluded from the Database"; a brand-new plain run with a single
#   space is inserted after it; the original "Scope" run
#   (w:rsidRPr="00951ACA") is left completely untouched.
# -----------------------------------------------------------------
$find2 = $d.Content
$find2.Find.ClearFormatting()
$found2 = $find2.Find.Execute("Outside ", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if ($found2) {
    $p2 = $find2.Paragraphs.Item(1)
    $find2.Text = ""

    $runs2 = '<w:r>' + $rPr + '<w:t>Excluded from the Database</w:t></w:r>' + `
             '<w:r>' + $rPr + '<w:t xml:space="preserve"> </w:t></w:r>'
    Insert-RunsAtParagraphStart $p2 $runs2
}

Write-Host "Edits applied."
